# Tiempos.docx — "Dos tiempos 10 y 100"
#
# 1) Resize the 4 table columns (tblGrid widths 2326/2326/2326/2326 dxa ->
#    1901/1985/2649/2769 dxa). Column widths on the Word object model are in
#    points, so dxa (twentieths of a point) / 20 = points.
# 2) Rename the two header cells to the "NoDirigido" test files.
# 3) Fill in the previously-empty timing cells: Prim/100 -> "1 ms",
#    Prim/1000 -> "15 ms", Kruskal/1000 -> "46 ms".
#
# The runtime re-seats the Table/Cell COM handles after any structural
# write, so each step below re-fetches $t (and the cell) from
# ActiveDocument.Tables before touching it, per the interop's own guidance.

$d = $word.ActiveDocument

# --- 1. Column widths (dxa / 20 = points) ---------------------------------
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 95.05   # 1901 dxa

$t = $d.Tables.Item(1)
$t.Columns.Item(2).Width = 99.25   # 1985 dxa

$t = $d.Tables.Item(1)
$t.Columns.Item(3).Width = 132.45  # 2649 dxa

$t = $d.Tables.Item(1)
$t.Columns.Item(4).Width = 138.45  # 2769 dxa

# --- 2. Header row: test-file names ---------------------------------------
$t = $d.Tables.Item(1)
$t.Cell(1, 3).Range.Text = "pruebaNoDirigido100.txt"

$t = $d.Tables.Item(1)
$t.Cell(1, 4).Range.Text = "pruebaNoDirigido1000.txt"

# --- 3. Prim row: fill in the 100 and 1000 timings -------------------------
$t = $d.Tables.Item(1)
$t.Cell(2, 2).Range.Text = "1 ms"

$t = $d.Tables.Item(1)
$t.Cell(2, 3).Range.Text = "15 ms"

# --- 4. Kruskal row: fill in the 1000 timing -------------------------------
$t = $d.Tables.Item(1)
$t.Cell(3, 3).Range.Text = "46 ms"

Write-Host "Tiempos.docx updated: columns resized, NoDirigido files named, timings filled in."
